$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.197608709335327
$ws.Range("B1").Value = 2.760985136032104
$ws.Range("C1").Value = 3.672125577926636
$ws.Range("D1").Value = 5.812647819519043
$ws.Range("E1").Value = 2.125693559646606
